$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI LR-pair values (Mdk-Sdc4) per commit "Natmi following Dr Hou advice".
# Recalculated with Ligand/Receptor-expressing cell counts (E, K) changed 1 -> 3,
# which cascades through avg/total expression, specificity, and edge-weight columns (G:T).
# Columns F and L (detection rate) remain 1, unchanged.

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 3
$row2[0,1] = 1
$row2[0,2] = 1.324023666666666
$row2[0,3] = 3.972071
$row2[0,4] = 0.01518042398701374
$row2[0,5] = 0.01518042398701374
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 3.438907666666667
$row2[0,9] = 10.316723
$row2[0,10] = 0.05825422340060618
$row2[0,11] = 0.05825422340060618
$row2[0,12] = 4.55319513814811
$row2[0,13] = 40.978756243333
$row2[0,14] = 0.0008843238102554191
$row2[0,15] = 0.0008843238102554191
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 3
$row3[0,1] = 1
$row3[0,2] = 1.324023666666666
$row3[0,3] = 3.972071
$row3[0,4] = 0.01518042398701374
$row3[0,5] = 0.01518042398701374
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 10.383857
$row3[0,9] = 31.151571
$row3[0,10] = 0.1758999031294962
$row3[0,11] = 0.1758999031294962
$row3[0,12] = 13.74847241928233
$row3[0,13] = 123.736251773541
$row3[0,14] = 0.002670235108780397
$row3[0,15] = 0.002670235108780397
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 3
$row4[0,1] = 1
$row4[0,2] = 1.324023666666666
$row4[0,3] = 3.972071
$row4[0,4] = 0.01518042398701374
$row4[0,5] = 0.01518042398701374
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 25.34077833333333
$row4[0,9] = 76.022335
$row4[0,10] = 0.4292663558501786
$row4[0,11] = 0.4292663558501786
$row4[0,12] = 33.55179024508722
$row4[0,13] = 301.966112205785
$row4[0,14] = 0.006516445285166026
$row4[0,15] = 0.006516445285166027
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 3
$row5[0,1] = 1
$row5[0,2] = 1.324023666666666
$row5[0,3] = 3.972071
$row5[0,4] = 0.01518042398701374
$row5[0,5] = 0.01518042398701374
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 19.86921833333334
$row5[0,9] = 59.60765500000001
$row5[0,10] = 0.336579517619719
$row5[0,11] = 0.336579517619719
$row5[0,12] = 26.30731531150056
$row5[0,13] = 236.765837803505
$row5[0,14] = 0.005109419782811896
$row5[0,15] = 0.005109419782811896
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 81.17653533333333
$row6[0,3] = 243.529606
$row6[0,4] = 0.9307191821270077
$row6[0,5] = 0.9307191821270075
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 3.438907666666667
$row6[0,9] = 10.316723
$row6[0,10] = 0.05825422340060618
$row6[0,11] = 0.05825422340060618
$row6[0,12] = 279.1586097112375
$row6[0,13] = 2512.427487401138
$row6[0,14] = 0.05421832315885618
$row6[0,15] = 0.05421832315885617
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 81.17653533333333
$row7[0,3] = 243.529606
$row7[0,4] = 0.9307191821270077
$row7[0,5] = 0.9307191821270075
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 10.383857
$row7[0,9] = 31.151571
$row7[0,10] = 0.1758999031294962
$row7[0,11] = 0.1758999031294962
$row7[0,12] = 842.9255346567808
$row7[0,13] = 7586.329811911026
$row7[0,14] = 0.1637134139769045
$row7[0,15] = 0.1637134139769045
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 81.17653533333333
$row8[0,3] = 243.529606
$row8[0,4] = 0.9307191821270077
$row8[0,5] = 0.9307191821270075
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 25.34077833333333
$row8[0,9] = 76.022335
$row8[0,10] = 0.4292663558501786
$row8[0,11] = 0.4292663558501786
$row8[0,12] = 2057.076587750001
$row8[0,13] = 18513.68928975001
$row8[0,14] = 0.3995264316315192
$row8[0,15] = 0.3995264316315192
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 81.17653533333333
$row9[0,3] = 243.529606
$row9[0,4] = 0.9307191821270077
$row9[0,5] = 0.9307191821270075
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 19.86921833333334
$row9[0,9] = 59.60765500000001
$row9[0,10] = 0.336579517619719
$row9[0,11] = 0.336579517619719
$row9[0,12] = 1612.914304081548
$row9[0,13] = 14516.22873673393
$row9[0,14] = 0.3132610133597276
$row9[0,15] = 0.3132610133597276
$ws.Range("E9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 1.192675
$row10[0,3] = 3.578025
$row10[0,4] = 0.0136744626508778
$row10[0,5] = 0.0136744626508778
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 3.438907666666667
$row10[0,9] = 10.316723
$row10[0,10] = 0.05825422340060618
$row10[0,11] = 0.05825422340060618
$row10[0,12] = 4.101499201341667
$row10[0,13] = 36.913492812075
$row10[0,14] = 0.0007965952021474809
$row10[0,15] = 0.0007965952021474808
$ws.Range("E10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 3
$row11[0,1] = 1
$row11[0,2] = 1.192675
$row11[0,3] = 3.578025
$row11[0,4] = 0.0136744626508778
$row11[0,5] = 0.0136744626508778
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 10.383857
$row11[0,9] = 31.151571
$row11[0,10] = 0.1758999031294962
$row11[0,11] = 0.1758999031294962
$row11[0,12] = 12.384566647475
$row11[0,13] = 111.461099827275
$row11[0,14] = 0.002405336655637319
$row11[0,15] = 0.002405336655637318
$ws.Range("E11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 3
$row12[0,1] = 1
$row12[0,2] = 1.192675
$row12[0,3] = 3.578025
$row12[0,4] = 0.0136744626508778
$row12[0,5] = 0.0136744626508778
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 25.34077833333333
$row12[0,9] = 76.022335
$row12[0,10] = 0.4292663558501786
$row12[0,11] = 0.4292663558501786
$row12[0,12] = 30.22331279870833
$row12[0,13] = 272.009815188375
$row12[0,14] = 0.005869986750351688
$row12[0,15] = 0.005869986750351688
$ws.Range("E12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 3
$row13[0,1] = 1
$row13[0,2] = 1.192675
$row13[0,3] = 3.578025
$row13[0,4] = 0.0136744626508778
$row13[0,5] = 0.0136744626508778
$row13[0,6] = 3
$row13[0,7] = 1
$row13[0,8] = 19.86921833333334
$row13[0,9] = 59.60765500000001
$row13[0,10] = 0.336579517619719
$row13[0,11] = 0.336579517619719
$row13[0,12] = 23.69751997570834
$row13[0,13] = 213.277679781375
$row13[0,14] = 0.004602544042741315
$row13[0,15] = 0.004602544042741315
$ws.Range("E13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,16
$row14[0,0] = 3
$row14[0,1] = 1
$row14[0,2] = 3.525915333333334
$row14[0,3] = 10.577746
$row14[0,4] = 0.04042593123510095
$row14[0,5] = 0.04042593123510094
$row14[0,6] = 3
$row14[0,7] = 1
$row14[0,8] = 3.438907666666667
$row14[0,9] = 10.316723
$row14[0,10] = 0.05825422340060618
$row14[0,11] = 0.05825422340060618
$row14[0,12] = 12.12529727181756
$row14[0,13] = 109.127675446358
$row14[0,14] = 0.002354981229347114
$row14[0,15] = 0.002354981229347114
$ws.Range("E14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,16
$row15[0,0] = 3
$row15[0,1] = 1
$row15[0,2] = 3.525915333333334
$row15[0,3] = 10.577746
$row15[0,4] = 0.04042593123510095
$row15[0,5] = 0.04042593123510094
$row15[0,6] = 3
$row15[0,7] = 1
$row15[0,8] = 10.383857
$row15[0,9] = 31.151571
$row15[0,10] = 0.1758999031294962
$row15[0,11] = 0.1758999031294962
$row15[0,12] = 36.61260061544068
$row15[0,13] = 329.513405538966
$row15[0,14] = 0.00711091738817393
$row15[0,15] = 0.007110917388173929
$ws.Range("E15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,16
$row16[0,0] = 3
$row16[0,1] = 1
$row16[0,2] = 3.525915333333334
$row16[0,3] = 10.577746
$row16[0,4] = 0.04042593123510095
$row16[0,5] = 0.04042593123510094
$row16[0,6] = 3
$row16[0,7] = 1
$row16[0,8] = 25.34077833333333
$row16[0,9] = 76.022335
$row16[0,10] = 0.4292663558501786
$row16[0,11] = 0.4292663558501786
$row16[0,12] = 89.34943888410112
$row16[0,13] = 804.1449499569101
$row16[0,14] = 0.01735349218314169
$row16[0,15] = 0.01735349218314169
$ws.Range("E16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,16
$row17[0,0] = 3
$row17[0,1] = 1
$row17[0,2] = 3.525915333333334
$row17[0,3] = 10.577746
$row17[0,4] = 0.04042593123510095
$row17[0,5] = 0.04042593123510094
$row17[0,6] = 3
$row17[0,7] = 1
$row17[0,8] = 19.86921833333334
$row17[0,9] = 59.60765500000001
$row17[0,10] = 0.336579517619719
$row17[0,11] = 0.336579517619719
$row17[0,12] = 70.05718158284779
$row17[0,13] = 630.5146342456302
$row17[0,14] = 0.01360654043443821
$row17[0,15] = 0.01360654043443821
$ws.Range("E17:T17").Value = $row17

